# Split the run-on "Programa" (PT/EN) and "Bibliografia" paragraphs into
# multiple <w:t> runs separated by <w:br/> line breaks, matching the
# target diff -- including which fragments keep xml:space="preserve".
#
# We replace each target paragraph's content with Range.InsertXML, driven
# by a Flat-OPC <pkg:package> payload wrapping the exact <w:p> we want.
# (Plain Find/Replace with the "^l" break code works for the text/break
# structure, but this host recomputes a whole run's xml:space from only
# its LAST <w:t> fragment, so it can't reproduce "earlier fragments keep
# their trailing-space preserve, the final fragment has none" -- which is
# exactly what this diff needs. InsertXML lets us state it directly.)

$d = $word.ActiveDocument

function Get-ParagraphByAnchor($anchorText) {
    foreach ($p in $d.Paragraphs) {
        if ($p.Range.Text.Contains($anchorText)) {
            return $p
        }
    }
    throw "No paragraph found containing: $anchorText"
}

function Set-ParagraphXml($anchorText, $pkgXml) {
    $p = Get-ParagraphByAnchor $anchorText
    [void]$p.Range.InsertXML($pkgXml)
}

$ptAnchor = 'A. MICROECONOMIA: 1. Introdução aos conceitos'
$ptPayload = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">A. MICROECONOMIA: 1. Introdução aos conceitos de Economia e fundamentos da análise microeconômica. 2. Teoria do consumidor e da demanda. 3. Teoria da firma e da oferta. 4. Custos e formação de preços. 5. Estruturas de Mercado 6. Comportamento estratégico e concorrência. 7. Tecnologia como fator de produção. 8. Sustentabilidade: recursos, custos e indicadores ambientais. </w:t><w:br/><w:t xml:space="preserve">B. MACROECONOMIA: 1. Fundamentos da análise macroeconômica. 2. Contabilidade nacional. 3. Equilíbrios clássicos e keynesiano. 4. Sistema monetário. 5. Política fiscal. 6. Economia mundial e comércio internacional. 7. Fundamentos da regressão como ferramenta para quantificar relações econômicas. 8. Setor público. </w:t><w:br/><w:t>C. DESENVOLVIMENTO ECONÔMICO: 1. Fatores de Crescimento. 2. Fontes de Desenvolvimento. 3. Financiamento do Desenvolvimento Econômico. 4. Um modelo de Crescimento Econômico. 5. O Processo de internacionalização e globalização.</w:t><w:br/><w:t>D. ECONOMIA INTERNACIONAL: 1. Fundamentos do Comércio Internacional. 2. Determinação das Taxas de Câmbio. 3. Políticas Externas. 4. Fatores determinantes do comportamento das importações e exportações.</w:t><w:br/><w:t>E. ECONOMIA BRASILEIRA: 1. A experiência histórica da industrialização brasileira. 2. A internacionalização da economia brasileira. 3. Teoria dos ciclos e realidade brasileira. 4. Os ciclos econômicos do Brasil ao longo de sua história recente.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
Set-ParagraphXml $ptAnchor $ptPayload

$enAnchor = 'A. MICROECONOMY: 1. Introduction to the concepts'
$enPayload = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:i/></w:rPr><w:t>A. MICROECONOMY: 1. Introduction to the concepts of Economics and fundamentals of microeconomic analysis. 2. Consumer and demand theory. 3. Firm and offer theory. 4. Costs and price formation. 5. Market Structures 6. Strategic behavior and competition. 7. Technology as a factor of production. 8. Sustainability: resources, costs and environmental indicators.</w:t><w:br/><w:t>B. MACROECONOMY: 1. Fundamentals of macroeconomic analysis. 2. National accounting. 3. Classical and Keynesian balances. 4. Monetary system. 5. Fiscal policy. 6. World economy and international trade. 7. Fundamentals of regression as a tool to quantify economic relationships. 8. Public sector.</w:t><w:br/><w:t>C. ECONOMIC DEVELOPMENT: 1. Growth factors. 2. Sources of Development. 3. Financing Economic Development. 4. A model of economic growth. 5. The internationalization and globalization process.</w:t><w:br/><w:t>D. INTERNATIONAL ECONOMY: 1. Fundamentals of International Trade. 2. Determination of Exchange Rates. 3. External policies. 4. Factors determining the behavior of imports and exports.</w:t><w:br/><w:t>E. BRAZILIAN ECONOMY: 1. The historical experience of Brazilian industrialization. 2. The internationalization of the Brazilian economy. 3. Cycle theory and Brazilian reality. 4. Brazil''s economic cycles throughout its recent history.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
Set-ParagraphXml $enAnchor $enPayload

$biblioAnchor = 'MANKIW, N.G. Introdução'
$biblioPayload = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>MANKIW, N.G. Introdução à economia. São Paulo: Thomson Learning, 2006.</w:t><w:br/><w:br/><w:t>SAMUELSON, P. Introdução à Economia. New York: Mc Graw-Hill Book Company.</w:t><w:br/><w:br/><w:t>BACHA, Edmar. Introdução à Macroeconomia: Uma perspectiva brasileira. Rio de Janeiro: Campus, 1987.</w:t><w:br/><w:br/><w:t>BACHA et al. Estado da Economia Mundial - Desafios e Respostas - Seminário em Homenagem a Pedro Malan. São Paulo: LTC, 2015.</w:t><w:br/><w:br/><w:t>FURTADO, C. Formação econômica do Brasil. São Paulo: Companhia Editora Nacional, 2003.</w:t><w:br/><w:br/><w:t>GREMAUD, A. P.; VASCONCELLOS, M. A. S.; TONETO JÚNIOR, R. Economia Brasileira Contemporânea. 8 ed. São Paulo: Atlas, 2017.</w:t><w:br/><w:br/><w:t>VASCONCELLOS, M. A. S.; GARCIA, M. E. Fundamentos de Economia. 6 ed. São Paulo: Saraiva, 2018.</w:t><w:br/><w:br/><w:t>VASCONCELLOS, M. A. S. ECONOMIA: Micro e Macro. São Paulo: Atlas, 2015.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
Set-ParagraphXml $biblioAnchor $biblioPayload

Write-Output "Done"
